$wb = $excel.ActiveWorkbook

# The "optimization_parameters" sheet had a stray row (row 16 - label "Sheet" with
# values 3 / 4) that was left over from earlier editing. Select the row and
# delete it entirely, which shifts every following row up by one.
$wsParams = $wb.Worksheets.Item("optimization_parameters")
$wsParams.Activate()
$wsParams.Rows.Item(16).Select()
$wsParams.Rows.Item(16).Delete()

# Finish on the "threshold_b" sheet (the last sheet touched/active when the
# file was saved), with its existing selection on A2.
$wsThreshold = $wb.Worksheets.Item("threshold_b")
$wsThreshold.Activate()
$wsThreshold.Range("A2").Select()
